$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("Xe vào bãi" bullet) - shape "Text 5"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shape6 = $s6.Shapes.Item(6)

# Move the textbox up slightly (xfrm offset y 2779276 -> 2680097 EMU)
$shape6.Top = 2680097 / 914400 * 72

$tr6 = $shape6.TextFrame.TextRange
$len6 = $tr6.Length
$whole6 = $tr6.Characters(1, $len6)
$whole6.Text = "Xe vào bãi: Kiểm tra khoảng cách, chọn ô trống, chụp ảnh xe, biển số, mở cổng và mô phỏng xe vào bến bằng động cơ DC-Motor"

# Re-apply (identical) formatting on sub-ranges so each becomes its own run,
# matching the run layout introduced by the edit.
$tr6.Characters(54,3).Font.Size = 18.5    # "ảnh"
$tr6.Characters(57,1).Font.Size = 18.5    # " "
$tr6.Characters(58,11).Font.Size = 18.5   # "xe, biển số"
$tr6.Characters(69,2).Font.Size = 18.5    # ", "
$tr6.Characters(71,2).Font.Size = 18.5    # "mở"
$tr6.Characters(73,1).Font.Size = 18.5    # " "
$tr6.Characters(74,44).Font.Size = 18.5   # "cổng và mô phỏng xe vào bến bằng động cơ DC-"
$tr6.Characters(118,5).Font.Size = 18.5   # "Motor"

# Add the new trailing (blank) paragraph containing a single space.
$null = $tr6.InsertAfter("`r ")

# ---------------------------------------------------------------------------
# Slide 7 ("Chụp ảnh xe vào/ra" bullet) - shape "Text 12"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(13)

$tr7 = $shape7.TextFrame.TextRange
$len7 = $tr7.Length
$whole7 = $tr7.Characters(1, $len7)
$whole7.Text = "Chụp ảnh xe vào/ra tự động để quản lý an ninh."

$tr7.Characters(17,2).Font.Size = 18.5   # "ra"
$tr7.Characters(19,8).Font.Size = 18.5   # " tự động"
$tr7.Characters(27,20).Font.Size = 18.5  # " để quản lý an ninh."
